# Commit: "update database and change read_price algorithm"
#
# The quarterly income-statement figures for this ticker are replaced with
# zeros (the read_price algorithm now reports "no data" as 0 / "-" instead
# of the previously scraped numbers), while row/column labels, styles and
# layout are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data columns D..M (quarterly periods) for every financial-statement row.
$dataCols = @(4,5,6,7,8,9,10,11,12,13)

# Rows whose entire D:M range becomes a plain numeric 0.
$zeroRows = @(11,12,13,14,16,17,18,19,20,22,24,25,26,27)
foreach ($r in $zeroRows) {
    foreach ($c in $dataCols) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# Rows whose entire D:M range becomes the text placeholder "-".
$dashRows = @(15,23)
foreach ($r in $dashRows) {
    foreach ($c in $dataCols) {
        $ws.Cells.Item($r, $c).Value = "-"
    }
}

# Row 21 is mixed: every quarter is zeroed except column G, which becomes "-".
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = "-"
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = 0
